$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the speaker links to point to individual bio pages instead of the
# generic Leadership_team.html page.
$ws.Range("F2").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/lauren_chenarides.html), [Drew Hanks](https://dataifa.github.io/difa-project/drew_hanks.html)"
$ws.Range("F4").Value = "[George Davis](https://dataifa.github.io/difa-project/george_davis.html), [Joe Cummins](https://www.josephrcummins.com/)"
$ws.Range("F6").Value = "[Drew Hanks](https://dataifa.github.io/difa-project/drew_hanks.html)"
$ws.Range("F7").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/lauren_chenarides.html)"
$ws.Range("F9").Value = "[Lauren Chenarides](https://dataifa.github.io/difa-project/lauren_chenarides.html), [Drew Hanks](https://dataifa.github.io/difa-project/drew_hanks.html)"

# Update the selected cell in the sheet view.
$ws.Range("H9").Select()
